$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1715
$ws1.Range("F10").Value = 116
$ws1.Range("F11").Value = 6032
$ws1.Range("F15").Value = 4826
$ws1.Range("F16").Value = 23
$ws1.Range("F21").Value = 48
$ws1.Range("F24").Value = 20
$ws1.Range("F25").Value = 3297
$ws1.Range("F26").Value = 137

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 42

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1715
$ws4.Range("F11").Value = 116
$ws4.Range("F12").Value = 6032
$ws4.Range("F16").Value = 4826
$ws4.Range("F17").Value = 23
$ws4.Range("F22").Value = 48
$ws4.Range("F25").Value = 20
$ws4.Range("F26").Value = 3297
$ws4.Range("F27").Value = 42
$ws4.Range("F28").Value = 137
